$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Comment" column header in R1, bold font + yellow fill to match
# the existing header style.
$ws.Range("R1").Value = "Comment"
$ws.Range("R1").Font.Bold = $true
$ws.Range("R1").Interior.Color = 65535

# Mark the "Run scripts" job rows (12-14) as archived.
$ws.Range("R12").Value = "Archive"
$ws.Range("R13").Value = "Archive"
$ws.Range("R14").Value = "Archive"
